# Fixing dollar sign add bug
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix existing salary values that had extra trailing '$' characters added by the bug.
# Temporarily mark the cell as Text so Excel doesn't reinterpret the
# "<number>$" value as a currency number, then restore the normal style so
# the cell keeps its original (unstyled) appearance.
$ws.Range('C2').NumberFormat = '@'
$ws.Range('C2').Value = '2300$'
$ws.Range('C2').Style = 'Normal'

$ws.Range('C3').NumberFormat = '@'
$ws.Range('C3').Value = '1000000$'
$ws.Range('C3').Style = 'Normal'

$ws.Range('C4').NumberFormat = '@'
$ws.Range('C4').Value = '3000$'
$ws.Range('C4').Style = 'Normal'

# Add a new row of test data demonstrating the fix
$ws.Range('A6').Value = 'TestingDollar'
$ws.Range('B6').Value = 'b''$2b$12$hPMP1PTGHNPEXPrmm112puz6ZbREw6wV9/cAfbDVr7rcS54VyerBi'''

$ws.Range('C6').NumberFormat = '@'
$ws.Range('C6').Value = '2750$'
$ws.Range('C6').Style = 'Normal'

$ws.Range('D6').Value = 'letstestthedollar@gmail.com'
